$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"25.73000000000058"
$ws.Range("H2").Value = [double]"0.009531053603840545"
$ws.Range("I2").Value = [double]"0.009531053603840545"
$ws.Range("L2").Value = [double]"38.33475478137952"
$ws.Range("M2").Value = "[4.809935402998946, 71.85957415976009]"
$ws.Range("N2").Value = [double]"0.02595148058490149"
$ws.Range("O2").Value = [double]"0.02595148058490149"
$ws.Range("P2").Value = [double]"2.06923720326535"
$ws.Range("Q2").Value = "[0.9622896416401172, 3.176184764890582]"
$ws.Range("R2").Value = [double]"0.0004805176601696548"
$ws.Range("S2").Value = [double]"0.0004805176601696548"
$ws.Range("T2").Value = [double]"69.95547342142032"
$ws.Range("U2").Value = "[52.39000525483122, 87.52094158800942]"
$ws.Range("V2").Value = [double]"3.209879029242302e-10"
$ws.Range("W2").Value = [double]"3.209879029242302e-10"
$ws.Range("X2").Value = [double]"17.25635635635675"
$ws.Range("Y2").Value = [double]"12.72334334334363"
$ws.Range("Z2").Value = [double]"21.78936936936987"
$ws.Range("F3").Value = [double]"25.73000000000058"
$ws.Range("H3").Value = [double]"0.0001340261053276448"
$ws.Range("I3").Value = [double]"0.0001340261053276448"
$ws.Range("L3").Value = [double]"60.93900876081536"
$ws.Range("M3").Value = "[26.689741369811827, 95.18827615181888]"
$ws.Range("N3").Value = [double]"0.0008287983821138045"
$ws.Range("O3").Value = [double]"0.0008287983821138045"
$ws.Range("Q3").Value = "[0.5597632555945786, 1.7673424137311944]"
$ws.Range("R3").Value = [double]"0.000336723297087671"
$ws.Range("S3").Value = [double]"0.000336723297087671"
$ws.Range("T3").Value = [double]"71.73828860675717"
$ws.Range("U3").Value = "[52.9853716476552, 90.49120556585913]"
$ws.Range("V3").Value = [double]"9.300586967242452e-10"
$ws.Range("W3").Value = [double]"9.300586967242452e-10"
$ws.Range("X3").Value = [double]"20.96518518518566"
$ws.Range("Y3").Value = [double]"18.49263263263305"
$ws.Range("Z3").Value = [double]"23.43773773773827"
$ws.Range("F4").Value = [double]"25.73000000000058"
$ws.Range("H4").Value = [double]"2.443748184066674e-05"
$ws.Range("I4").Value = [double]"2.443748184066674e-05"
$ws.Range("L4").Value = [double]"60.8060366848594"
$ws.Range("M4").Value = "[30.64485501536278, 90.96721835435602]"
$ws.Range("N4").Value = [double]"0.0001931365006131447"
$ws.Range("O4").Value = [double]"0.0001931365006131447"
$ws.Range("P4").Value = [double]"0.5597632555945786"
$ws.Range("Q4").Value = "[0.04402632347373281, 1.0755001877154244]"
$ws.Range("R4").Value = [double]"0.03405548766870359"
$ws.Range("S4").Value = [double]"0.03405548766870359"
$ws.Range("T4").Value = [double]"67.30447785844748"
$ws.Range("U4").Value = "[51.18928199814385, 83.41967371875111]"
$ws.Range("V4").Value = [double]"8.738343382219682e-11"
$ws.Range("W4").Value = [double]"8.738343382219682e-11"
$ws.Range("X4").Value = [double]"23.43773773773827"
$ws.Range("Y4").Value = [double]"21.32576576576625"
$ws.Range("Z4").Value = [double]"25.54970970971028"
$ws.Range("F5").Value = [double]"25.73000000000058"
$ws.Range("H5").Value = [double]"0.0007855926962417881"
$ws.Range("I5").Value = [double]"0.0007855926962417881"
$ws.Range("L5").Value = [double]"44.00623581065069"
$ws.Range("M5").Value = "[16.543352713407685, 71.46911890789369]"
$ws.Range("N5").Value = [double]"0.002332525280865916"
$ws.Range("O5").Value = [double]"0.002332525280865916"
$ws.Range("P5").Value = [double]"0.3962369112635784"
$ws.Range("Q5").Value = "[-0.2578684660604216, 1.0503422885875784]"
$ws.Range("R5").Value = [double]"0.2287894609238514"
$ws.Range("S5").Value = [double]"0.2287894609238514"
$ws.Range("T5").Value = [double]"72.41851830346411"
$ws.Range("U5").Value = "[57.67008223108063, 87.1669543758476]"
$ws.Range("V5").Value = [double]"7.351896869067787e-13"
$ws.Range("W5").Value = [double]"7.351896869067787e-13"
$ws.Range("X5").Value = [double]"24.10738738738793"
$ws.Range("Y5").Value = [double]"21.42878878878928"
$ws.Range("Z5").Value = [double]"26.78598598598658"
$ws.Range("F6").Value = [double]"25.73000000000058"
$ws.Range("H6").Value = [double]"5.64921400325602e-07"
$ws.Range("I6").Value = [double]"5.64921400325602e-07"
$ws.Range("L6").Value = [double]"52.42286068724266"
$ws.Range("M6").Value = "[33.78915631186774, 71.05656506261758]"
$ws.Range("N6").Value = [double]"9.741735369495075e-07"
$ws.Range("O6").Value = [double]"9.741735369495075e-07"
$ws.Range("P6").Value = [double]"-0.1257894956392311"
$ws.Range("Q6").Value = "[-0.5408948312486928, 0.28931583997023047]"
$ws.Range("R6").Value = [double]"0.5447105662075939"
$ws.Range("S6").Value = [double]"0.5447105662075939"
$ws.Range("T6").Value = [double]"53.50040053676737"
$ws.Range("U6").Value = "[42.01919552042839, 64.98160555310633]"
$ws.Range("V6").Value = [double]"3.649969215757665e-12"
$ws.Range("W6").Value = [double]"3.649969215757665e-12"
$ws.Range("X6").Value = [double]"0.5151151151151296"
$ws.Range("Y6").Value = [double]"-1.184764764764786"
$ws.Range("Z6").Value = [double]"2.214994994995045"
$ws.Range("B7").Value = [double]"0"
$ws.Range("F7").Value = [double]"25.73000000000058"
$ws.Range("H7").Value = [double]"0.001471509534348936"
$ws.Range("I7").Value = [double]"0.001471509534348936"
$ws.Range("L7").Value = [double]"37.54277129114901"
$ws.Range("M7").Value = "[13.720926365181484, 61.364616217116534]"
$ws.Range("N7").Value = [double]"0.00270995785604744"
$ws.Range("O7").Value = [double]"0.00270995785604744"
$ws.Range("P7").Value = [double]"-0.2138421425866923"
$ws.Range("Q7").Value = "[-0.9434212172942322, 0.5157369321208476]"
$ws.Range("R7").Value = [double]"0.5579154612048596"
$ws.Range("S7").Value = [double]"0.5579154612048596"
$ws.Range("T7").Value = [double]"64.06274753576611"
$ws.Range("U7").Value = "[50.58290684276545, 77.54258822876676]"
$ws.Range("V7").Value = [double]"2.010169808386308e-12"
$ws.Range("W7").Value = [double]"2.010169808386308e-12"
$ws.Range("X7").Value = [double]"0.8756956956957147"
$ws.Range("Y7").Value = [double]"-2.111971971972022"
$ws.Range("Z7").Value = [double]"3.863363363363452"
$ws.Range("F8").Value = [double]"25.73000000000058"
$ws.Range("H8").Value = [double]"2.200730461088263e-06"
$ws.Range("I8").Value = [double]"2.200730461088263e-06"
$ws.Range("L8").Value = [double]"50.6159517607102"
$ws.Range("M8").Value = "[26.256630146215613, 74.9752733752048]"
$ws.Range("N8").Value = [double]"0.0001305058360372424"
$ws.Range("O8").Value = [double]"0.0001305058360372424"
$ws.Range("P8").Value = [double]"-0.8931054190385392"
$ws.Range("Q8").Value = "[-1.3333686537758478, -0.4528421843012307]"
$ws.Range("R8").Value = [double]"0.000178449511908374"
$ws.Range("S8").Value = [double]"0.000178449511908374"
$ws.Range("T8").Value = [double]"59.36882684329316"
$ws.Range("U8").Value = "[46.882452669686685, 71.85520101689963]"
$ws.Range("V8").Value = [double]"1.981970143560829e-12"
$ws.Range("W8").Value = [double]"1.981970143560829e-12"
$ws.Range("X8").Value = [double]"3.6573173173174"
$ws.Range("Y8").Value = [double]"1.854414414414455"
$ws.Range("Z8").Value = [double]"5.460220220220345"
$ws.Range("F9").Value = [double]"24.5800000000004"
$ws.Range("H9").Value = [double]"0.000363289543778289"
$ws.Range("I9").Value = [double]"0.000363289543778289"
$ws.Range("L9").Value = [double]"54.19218480064631"
$ws.Range("M9").Value = "[23.71478046739088, 84.66958913390174]"
$ws.Range("N9").Value = [double]"0.000834620176510148"
$ws.Range("O9").Value = [double]"0.000834620176510148"
$ws.Range("P9").Value = [double]"-1.685579241565694"
$ws.Range("Q9").Value = "[-2.3648425180175416, -1.006315965113847]"
$ws.Range("R9").Value = [double]"9.261790392178781e-06"
$ws.Range("S9").Value = [double]"9.261790392178781e-06"
$ws.Range("T9").Value = [double]"66.37962372002822"
$ws.Range("U9").Value = "[48.49258107683802, 84.26666636321842]"
$ws.Range("V9").Value = [double]"2.027569223628234e-09"
$ws.Range("W9").Value = [double]"2.027569223628234e-09"
$ws.Range("X9").Value = [double]"6.594034034034141"
$ws.Range("Y9").Value = [double]"3.9367367367368"
$ws.Range("Z9").Value = [double]"9.251331331331484"
$ws.Range("F10").Value = [double]"24.5800000000004"
$ws.Range("H10").Value = [double]"6.342040106632751e-06"
$ws.Range("I10").Value = [double]"6.342040106632751e-06"
$ws.Range("L10").Value = [double]"55.89052320668113"
$ws.Range("M10").Value = "[31.938395800730305, 79.84265061263196]"
$ws.Range("N10").Value = [double]"2.481922391028846e-05"
$ws.Range("O10").Value = [double]"2.481922391028846e-05"
$ws.Range("P10").Value = [double]"-1.42142130072331"
$ws.Range("Q10").Value = "[-1.9497371824080805, -0.8931054190385392]"
$ws.Range("R10").Value = [double]"2.254533585599461e-06"
$ws.Range("S10").Value = [double]"2.254533585599461e-06"
$ws.Range("T10").Value = [double]"59.38699400258037"
$ws.Range("U10").Value = "[45.060432866458655, 73.71355513870209]"
$ws.Range("V10").Value = [double]"1.076465583338404e-10"
$ws.Range("W10").Value = [double]"1.076465583338404e-10"
$ws.Range("X10").Value = [double]"5.560640640640731"
$ws.Range("Y10").Value = [double]"3.493853853853908"
$ws.Range("Z10").Value = [double]"7.627427427427554"
$ws.Range("F11").Value = [double]"24.5800000000004"
$ws.Range("H11").Value = [double]"0.0002734288872463519"
$ws.Range("I11").Value = [double]"0.0002734288872463519"
$ws.Range("L11").Value = [double]"48.41697067194163"
$ws.Range("M11").Value = "[19.972855473379852, 76.86108587050342]"
$ws.Range("N11").Value = [double]"0.001309224229875339"
$ws.Range("O11").Value = [double]"0.001309224229875339"
$ws.Range("P11").Value = [double]"-1.974895081535926"
$ws.Range("Q11").Value = "[-2.616421509296004, -1.3333686537758473]"
$ws.Range("R11").Value = [double]"1.57370266595791e-07"
$ws.Range("S11").Value = [double]"1.57370266595791e-07"
$ws.Range("T11").Value = [double]"62.46160690564658"
$ws.Range("U11").Value = "[46.8905527737807, 78.03266103751247]"
$ws.Range("V11").Value = [double]"2.642888130566234e-10"
$ws.Range("W11").Value = [double]"2.642888130566234e-10"
$ws.Range("X11").Value = [double]"7.725845845845971"
$ws.Range("Y11").Value = [double]"5.216176176176258"
$ws.Range("Z11").Value = [double]"10.23551551551568"
